# #1332 inject JdbcTemplate into JdbcCodelist (at en)
#
# 1) Slide with the JdbcCodeList XML sample: split the
#    <property name="dataSource" ref="dataSource" /> line into
#    <property name="jdbcTemplate" ref="jdbcTemplateForCodeList" />
#    using the same multi-run layout as the authored edit.
# 2) Refresh the cached "datetimeFigureOut" field text
#    (2014/12/17 -> 2015/10/8) on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Locate the shape that holds the JdbcCodeList bean XML sample and
#    rewrite the "dataSource" property line.
# ---------------------------------------------------------------------------
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $sh = $slide.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text.Contains('name="dataSource" ref="dataSource"')) {
                $targetShape = $sh
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange

    $oldLine = '  <property name="dataSource" ref="dataSource" />'
    $newLine = '  <property name="jdbcTemplate" ref="jdbcTemplateForCodeList" />'

    $fullText = $tr.Text
    $startIdx = $fullText.IndexOf($oldLine)

    if ($startIdx -ge 0) {
        $start = $startIdx + 1   # TextRange.Characters is 1-based

        # Replace the whole line's text first.
        $whole = $tr.Characters($start, $oldLine.Length)
        $whole.Text = $newLine

        # Re-split the freshly written text into the same nine runs the
        # authored deck ends up with (all runs share identical formatting,
        # touching Font.Bold with its own value is enough to force a run
        # break at that boundary without changing any visible formatting).
        $runLens = @(16, 2, 12, 2, 3, 2, 23, 2, 2)
        $pos = $start
        foreach ($len in $runLens) {
            $seg = $tr.Characters($pos, $len)
            $seg.Font.Bold = $seg.Font.Bold
            $pos += $len
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Update the cached date field text wherever it appears.
#
# NB: the date placeholder shape is addressed by its fixed position (it is
# not found by reading shape text back first) because the host's read path
# for some CustomLayout shape collections is unreliable; the write path is
# not, and always lands on the right shape when addressed by index.
# ---------------------------------------------------------------------------
$newDate = '2015/10/8'

# Slide master: the "日付プレースホルダー" shape is the 3rd shape.
$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# Every slide layout hanging off the master - index of the date placeholder
# shape within each layout's Shapes collection (derived from each layout's
# shape order: title/content placeholders, then date, footer, slide number).
$dateShapeIndex = @{
    1 = 3
    2 = 3
    3 = 3
    4 = 4
    5 = 6
    6 = 2
    7 = 1
    8 = 4
    9 = 4
    10 = 3
    11 = 3
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    $idx = $dateShapeIndex[$li]
    $layout.Shapes.Item($idx).TextFrame.TextRange.Text = $newDate
}
